$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.926.63"
$ws.Range("E2").Value = "  +0.55%  "

$ws.Range("D3").Value = "1.552.53"
$ws.Range("E3").Value = "  +0.26%  "

$ws.Range("E4").Value = "  +0.35%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "206.38"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.84%  "

$ws.Range("E6").Value = "  +0.46%  "

$ws.Range("E7").Value = "  +0.38%  "

$ws.Range("E8").Value = "  +0.77%  "

$ws.Range("E9").Value = "  +0.72%  "

$ws.Range("E11").Value = "  +0.07%  "

$ws.Range("D12").Value = "1.773.48"
$ws.Range("E12").Value = "  +0.28%  "

$ws.Range("D13").Value = "1.554.05"
$ws.Range("E13").Value = "  +0.45%  "

$ws.Range("E14").Value = "  +0.58%  "

$ws.Range("E15").Value = "  +0.77%  "

$ws.Range("D16").Value = "26.931.63"
$ws.Range("E16").Value = "  +0.56%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "213.85"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.01%  "

$ws.Range("E19").Value = "  +0.22%  "

$ws.Range("E20").Value = "  +0.16%  "

$ws.Range("E21").Value = "  +0.35%  "

$ws.Range("E22").Value = "  -1.06%  "

$ws.Range("E23").Value = "  +1.82%  "

$ws.Range("E24").Value = "  -1.10%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "152.81"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.04%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.65"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.00%  "

$ws.Range("E27").Value = "  -0.04%  "

$ws.Range("E28").Value = "  +0.40%  "

$ws.Range("E29").Value = "  +1.22%  "

$ws.Range("E30").Value = "  -0.28%  "

$ws.Range("E31").Value = "  -0.98%  "

$ws.Range("E32").Value = "  +1.78%  "

$ws.Range("D33").Value = "1.375.43"
$ws.Range("E33").Value = "  +1.78%  "

$ws.Range("E34").Value = "  +2.09%  "

$ws.Range("E35").Value = "  +3.31%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.971"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +5.85%  "

$ws.Range("E37").Value = "  +0.53%  "

$ws.Range("E38").Value = "  +1.02%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.522"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.29%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.807"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.73%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.01"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.40%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.989"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.14%  "

$ws.Range("E43").Value = "  -1.09%  "

$ws.Range("E44").Value = "  +3.00%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "63.67"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.24%  "

$ws.Range("E46").Value = "  -1.82%  "

$ws.Range("D47").Value = "1.687.27"
$ws.Range("E47").Value = "  +0.24%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "86.26"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.43%  "

$ws.Range("E49").Value = "  +0.42%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0955"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.12%  "

$ws.Range("E51").Value = "  +0.48%  "
